# Append a new job row at the top of the data (row 9) on the "ランサーズ" sheet,
# pushing the previous rows 9-13 down to 10-14, and refresh the "取得日時"
# timestamp on every data row to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('ランサーズ')

$newTimestamp = '2026-01-08 01:27:20'

# --- 1. Shift existing rows 9..13 down to 10..14 (bottom-up to avoid overwrite) ---
for ($r = 13; $r -ge 9; $r--) {
    $dst = $r + 1
    $ws.Range("A$dst").Value2 = $ws.Range("A$r").Value2
    $ws.Range("B$dst").Value2 = $ws.Range("B$r").Value2
    $ws.Range("C$dst").Value2 = $ws.Range("C$r").Value2
    $ws.Range("D$dst").Value2 = $ws.Range("D$r").Value2
    $ws.Range("E$dst").Value2 = $ws.Range("E$r").Value2
    $ws.Range("F$dst").Value2 = $ws.Range("F$r").Value2
    $ws.Range("G$dst").Value2 = $ws.Range("G$r").Value2
    $ws.Range("H$dst").Value2 = $ws.Range("H$r").Value2
}

# --- 2. Write the brand-new row 9 ---
$ws.Range('A9').Value = $newTimestamp
$ws.Range('B9').Value = '進行管理およびチームディレクションを担当'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '~ 5,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('F9').Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Range('G9').Value = 30
$ws.Range('H9').Value = '◇管理'

# --- 3. Refresh the timestamp column for every data row (2..14) ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("A$r").Value = $newTimestamp
}

# --- 4. Rebuild the URL hyperlinks for F2..F14 so they line up with the shifted data ---
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Range("F$r")
    $url = $cell.Value2
    $ws.Hyperlinks.Add($cell, $url)
}
